$wb = $excel.ActiveWorkbook

# Helper-free approach: force numeric-looking labels ("2023", "2022", ...)
# to be written as TEXT (matching every other year label already in the
# sheet), not auto-converted to numbers, while leaving no stray number
# format / style behind.

# --- Sheet 1: "Data" -------------------------------------------------------
# Insert two new rows right after the header row for the years 2023 and
# 2022; the previously-existing historical rows shift down unchanged.
$wsData = $wb.Worksheets.Item("Data")

$wsData.Rows("2:3").Insert()

$wsData.Cells.Item(2,1).NumberFormat = "@"
$wsData.Cells.Item(2,1).Value = "2023"
$wsData.Cells.Item(2,1).Style = "Normal"
$wsData.Cells.Item(2,2).Value = 4.1

$wsData.Cells.Item(3,1).NumberFormat = "@"
$wsData.Cells.Item(3,1).Value = "2022"
$wsData.Cells.Item(3,1).Style = "Normal"
$wsData.Cells.Item(3,2).Value = 4.1

# --- Sheet 2: "Metadata" ---------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# A1 used to hold a truly-empty string; it now matches B1 (" ").
$wsMeta.Range("A1").Value = " "

# Insert a new "actualizacion" row right before the "cita" row (row 9),
# pushing "cita" and the closing row down by one.
$wsMeta.Rows("9:9").Insert()

$wsMeta.Cells.Item(9,1).Value = "actualizacion"
$wsMeta.Cells.Item(9,2).Value = "Julio 2025"
